# CWL: add "restore character data" UI prompt strings (CN/JP) and populate
# worksheet rows 115-117 with the new localisation keys.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 115: cwl_ui_chara_restore ------------------------------------------------
$ws.Range("A115").Value = "cwl_ui_chara_restore"

# --- Row 116: cwl_ui_chara_restore_yes --------------------------------------------
$ws.Range("A116").Value = "cwl_ui_chara_restore_yes"

# --- Row 117: cwl_ui_chara_restore_no ---------------------------------------------
$ws.Range("A117").Value = "cwl_ui_chara_restore_no"

# --- D116 / D117: short CN strings --------------------------------------------
$ws.Range("D116").Value = "恢复角色数据"
$ws.Range("D117").Value = "保持现状"

# --- D115: CN rich text (two runs: 宋体 body + Cascadia Code placeholder) ---------
$ws.Range("D115").Value = "下列角色数据可供恢复:`n{0}"

# --- C115: JP rich text (two runs: Yu Gothic body + Cascadia Code placeholder) ----
$ws.Range("C115").Value = "次のキャラクターのデータが復元可能です：`n{0}"

# --- C116 / C117: short JP strings ----------------------------------------------
$ws.Range("C116").Value = "キャラクターデータを復元"
$ws.Range("C117").Value = "このままにする"

# --- Formatting -------------------------------------------------------------
# D115: wraps (matches existing "C112/D80"-style wrapped Cascadia Code cells)
$ws.Range("D115").Font.Name = "Cascadia Code"
$ws.Range("D115").Font.Size = 15.8
$ws.Range("D115").VerticalAlignment = -4108
$ws.Range("D115").WrapText = $true

# C115: new style - Cascadia Code font, wrapped, vertical centered
$ws.Range("C115").Font.Name = "Cascadia Code"
$ws.Range("C115").Font.Size = 15.8
$ws.Range("C115").VerticalAlignment = -4108
$ws.Range("C115").WrapText = $true

# D116/D117: 宋体, vertical centered, no wrap
$ws.Range("D116").Font.Name = "宋体"
$ws.Range("D116").Font.Size = 15.8
$ws.Range("D116").VerticalAlignment = -4108
$ws.Range("D117").Font.Name = "宋体"
$ws.Range("D117").Font.Size = 15.8
$ws.Range("D117").VerticalAlignment = -4108

# C116/C117: Yu Gothic, vertical centered, no wrap
$ws.Range("C116").Font.Name = "Yu Gothic"
$ws.Range("C116").Font.Size = 15.8
$ws.Range("C116").VerticalAlignment = -4108
$ws.Range("C117").Font.Name = "Yu Gothic"
$ws.Range("C117").Font.Size = 15.8
$ws.Range("C117").VerticalAlignment = -4108

# --- Rich-text runs for the two mixed-font cells -----------------------------
# D115: "下列角色数据可供恢复:" (宋体) + "\n{0}" (Cascadia Code)
$r1 = $ws.Range("D115").Characters(1, 11)
$r1.Font.Name = "宋体"
$r1.Font.Size = 15.8
$r1.Font.ColorIndex = -4105
$r2 = $ws.Range("D115").Characters(12, 4)
$r2.Font.Name = "Cascadia Code"
$r2.Font.Size = 15.8
$r2.Font.ColorIndex = -4105

# C115: "次のキャラクターのデータが復元可能です：" (Yu Gothic) + "\n{0}" (Cascadia Code)
$r3 = $ws.Range("C115").Characters(1, 20)
$r3.Font.Name = "Yu Gothic"
$r3.Font.Size = 15.8
$r3.Font.ColorIndex = -4105
$r4 = $ws.Range("C115").Characters(21, 4)
$r4.Font.Name = "Cascadia Code"
$r4.Font.Size = 15.8
$r4.Font.ColorIndex = -4105

# --- Row heights now reflect the taller wrapped/CJK content -----------------
$ws.Rows(115).RowHeight = 48.75
$ws.Rows(116).RowHeight = 25.5
$ws.Rows(117).RowHeight = 25.5

# --- View state: scroll position + selection now sit on the new rows --------
$ws.Application.ActiveWindow.ScrollRow = 87
$ws.Range("C115:C117").Select() | Out-Null
Write-Host "Done writing values"
